# Insert a new weekly price record at row 51 for Repollo (Crespo record),
# Terminal Hortofrutícola Agro Chillán, shifting all subsequent rows down
# by one (the sheet's last row moves from 188 to 189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 51 - everything below shifts down.
$ws.Rows.Item(51).EntireRow.Insert()

# Populate the newly inserted row 51 with the new record's data.
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 44623
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = 100112006
$ws.Cells.Item(51, 7).Value = "Repollo"
$ws.Cells.Item(51, 8).Value = "Crespo record"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 200
$ws.Cells.Item(51, 11).Value = 850
$ws.Cells.Item(51, 12).Value = 900
$ws.Cells.Item(51, 13).Value = 875
$ws.Cells.Item(51, 14).Value = "$/unidad"
$ws.Cells.Item(51, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(51, 16).Value = 875
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
